# Commit: "added DB3 and info_PLC"
# - set DB2 (sheet1) selection to A7:C7
# - set DB1 (sheet2) selection to A1:E1, add pageSetup (paper size 9 / portrait)
# - add new worksheet "DB3" (PLC data-block listing) after DB1
# - add new worksheet "info_PLC" (connection info) after DB3, ends up active/tabSelected

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# DB2 (sheet1): change selection from A1:E7 to A7:C7
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("DB2")
$ws1.Activate()
$ws1.Range("A7:C7").Select()

# ---------------------------------------------------------------------------
# DB1 (sheet2): change selection from D5 to A1:E1, set page setup
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("DB1")
$ws2.Activate()
$ws2.Range("A1:E1").Select()
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Add DB3 sheet at the end (after DB1)
# ---------------------------------------------------------------------------
$count = $wb.Worksheets.Count
$db3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($count))
$db3.Name = "DB3"

# numeric value first, then format the range as text so the literal
# number 2 is kept as a number (matching <c r="D2" s="2"><v>2</v></c>)
$db3.Range("D2").Value2 = 2
$db3.Range("A2:E12").NumberFormat = "@"
$db3.Range("B1:E1").NumberFormat = "@"

$db3.Range("A1").Value = "Adress"
$db3.Range("B1").Value = "Name"
$db3.Range("C1").Value = "Type"
$db3.Range("D1").Value = "Initial value"
$db3.Range("E1").Value = "Comment"

$db3.Range("A2").Value = "0.0"
$db3.Range("B2").Value = "DB_VAR"
$db3.Range("C2").Value = "INT"
$db3.Range("E2").Value = "Temporary placeholder variable"

$db3.Range("A3").Value = "2.0"
$db3.Range("A4").Value = "6.0"
$db3.Range("A5").Value = "8.0"

# Name/Type/Initial-value columns entered row-by-row (matches the original
# authoring order, which determines shared-string table ordering)
$db3.Range("B3").Value = "DB_DWORD"
$db3.Range("C3").Value = "DWORD"
$db3.Range("D3").Value = "DW#16#57"

$db3.Range("B4").Value = "DB_Byte"
$db3.Range("C4").Value = "BYTE"
$db3.Range("D4").Value = "B#16#39"

$db3.Range("B5").Value = "DB_WORD"
$db3.Range("C5").Value = "WORD"
$db3.Range("D5").Value = "W#16#150"

$db3.Range("B6").Value = "DB_DINT"
$db3.Range("C6").Value = "DINT"
$db3.Range("D6").Value = "L#300"

$db3.Range("B7").Value = "DB_S5TIME"
$db3.Range("C7").Value = "S5TIME"
$db3.Range("D7").Value = "S5T#20MS"

$db3.Range("B8").Value = "DB_TIME"
$db3.Range("C8").Value = "TIME"
$db3.Range("D8").Value = "T#14MS"

$db3.Range("B9").Value = "DB_DATE"
$db3.Range("C9").Value = "DATE"
$db3.Range("D9").Value = "D#2018-5-30"

$db3.Range("B10").Value = "DB_TOD"
$db3.Range("C10").Value = "TIME_OF_DAY"
$db3.Range("D10").Value = "TOD#12:30:15.000"

$db3.Range("B11").Value = "DB_CHAR"
$db3.Range("C11").Value = "CHAR"
# leading apostrophe is Excel's "force text" quote prefix; double it so the
# literal value keeps its own leading quote: "'B'"
$db3.Range("D11").Value = "''B'"

# Adress column for rows 6-11 entered afterwards (new shared strings 59-64)
$db3.Range("A6").Value = "10.0"
$db3.Range("A7").Value = "14.0"
$db3.Range("A8").Value = "16.0"
$db3.Range("A9").Value = "20.0"
$db3.Range("A10").Value = "22.0"
$db3.Range("A11").Value = "26.0"

$db3.Range("A12").Value = "264.0"
$db3.Range("C12").Value = "END_STRUCT"

$db3.Range("A12").Select()

# ---------------------------------------------------------------------------
# Add info_PLC sheet at the end (after DB3) - becomes the active sheet
# ---------------------------------------------------------------------------
$count = $wb.Worksheets.Count
$info = $wb.Worksheets.Add($null, $wb.Worksheets.Item($count))
$info.Name = "info_PLC"

$info.Range("B2").Value2 = 0
$info.Range("B3").Value2 = 2
$info.Range("B1:B3").NumberFormat = "@"

$info.Range("A1").Value = "IP_adress"
$info.Range("B1").Value = "10.32.0.95"
$info.Range("A2").Value = "rack"
$info.Range("A3").Value = "slot"

$info.Range("D2").Select()
